$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 117842.3560886971
$ws.Range("E2").Value = 0.03480887865753751
$ws.Range("F2").Value = 0.2083876515397173
$ws.Range("G2").Value = -1.517125616962864
$ws.Range("H2").Value = 13.01025193176428
$ws.Range("D4").Value = 119413.9730069361
$ws.Range("E4").Value = 0.003636777061898183
$ws.Range("F4").Value = 0.2581598333271358
$ws.Range("G4").Value = -1.594843558567075
$ws.Range("H4").Value = 14.6061876952162
$ws.Range("D5").Value = 120188.0826928579
$ws.Range("E5").Value = [double]"-4.596591229468189e-06"
$ws.Range("F5").Value = 0.239418996831256
$ws.Range("G5").Value = -0.6897795567840518
$ws.Range("H5").Value = 7.776512527321323
$ws.Range("D6").Value = 120679.5346051632
$ws.Range("E6").Value = -0.0129212589388291
$ws.Range("F6").Value = 0.2686272240569052
$ws.Range("G6").Value = -1.113558203396932
$ws.Range("H6").Value = 9.766889088260266
$ws.Range("D7").Value = 122369.2132783464
$ws.Range("E7").Value = -0.03295327730884902
$ws.Range("F7").Value = 0.2243752785052507
$ws.Range("G7").Value = -0.9281486600866998
$ws.Range("H7").Value = 7.708840126100794
$ws.Range("D8").Value = 123794.9935553745
$ws.Range("E8").Value = -0.06290132379193823
$ws.Range("F8").Value = 0.3341681536268916
$ws.Range("G8").Value = -1.569242060850352
$ws.Range("H8").Value = 10.01434846283351
$ws.Range("D9").Value = 125038.2742885055
$ws.Range("E9").Value = -0.1053904517232215
$ws.Range("F9").Value = 0.4477531760926195
$ws.Range("G9").Value = -1.94453028825256
$ws.Range("H9").Value = 9.866577669656291
$ws.Range("D11").Value = 117836.9972784511
$ws.Range("E11").Value = 0.216211654142799
$ws.Range("F11").Value = 0.1244676369703484
$ws.Range("G11").Value = -1.516108179991855
$ws.Range("H11").Value = 12.48148922470957
$ws.Range("D12").Value = 117921.4858758598
$ws.Range("E12").Value = 0.1112488352452739
$ws.Range("F12").Value = 0.1752522630527318
$ws.Range("G12").Value = -0.830784814255098
$ws.Range("H12").Value = 7.917398500880876
$ws.Range("D13").Value = 117841.216714712
$ws.Range("E13").Value = 0.1847764920115214
$ws.Range("F13").Value = 0.1470053744034521
$ws.Range("G13").Value = -1.149054835375477
$ws.Range("H13").Value = 9.133740400557254
$ws.Range("D15").Value = 117891.9396738128
$ws.Range("E15").Value = 0.1482616109569921
$ws.Range("F15").Value = 0.1888558532581152
$ws.Range("G15").Value = -1.029176513388569
$ws.Range("H15").Value = 9.342419389652465
$ws.Range("D16").Value = 117775.0211533045
$ws.Range("E16").Value = 0.1278979360488205
$ws.Range("F16").Value = 0.199233787984637
$ws.Range("G16").Value = -1.887144836242889
$ws.Range("H16").Value = 15.51763069064654
$ws.Range("D18").Value = 117946.8995725916
$ws.Range("E18").Value = 0.06722534014136954
$ws.Range("F18").Value = 0.1744646347608071
$ws.Range("G18").Value = -0.643348862800846
$ws.Range("H18").Value = 7.01610172852216
$ws.Range("D19").Value = 118144.0492536415
$ws.Range("E19").Value = 0.05651701431838177
$ws.Range("F19").Value = 0.1781844537369143
$ws.Range("G19").Value = -0.52351876643772
$ws.Range("H19").Value = 6.422084728475692
